{"js": "// Add the new chapter title paragraph (\"General Discussion, Conclusions, and\n// Future Perspectives\") plus a following blank paragraph, right before the\n// existing first paragraph (\"Note: some generic words ...\").\n//\n// We build the new content with raw OOXML (flat-OPC wrapped, as required by\n// Office.js's insertOoxml) so that the run boundaries exactly match the\n// source document (the title text is split across several runs, most likely\n// left over from Word's auto-capitalize-first-letter-of-sentence typing\n// history: \"General Discussion, \" | \"C\" | \"onclusions, and \" | \"F\" |\n// \"uture \" | \"P\" | \"erspectives\"), and the spacer paragraph is a fully empty\n// <w:p/> (no run), matching the target diff precisely.\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">General Discussion, </w:t></w:r>' +\n  '<w:r><w:t>C</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">onclusions, and </w:t></w:r>' +\n  '<w:r><w:t>F</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">uture </w:t></w:r>' +\n  '<w:r><w:t>P</w:t></w:r>' +\n  '<w:r><w:t>erspectives</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst body = context.document.body;\nconst startRange = body.getRange(Word.RangeLocation.start);\nstartRange.insertOoxml(flatOpcXml, Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Add the new chapter title paragraph (\"General Discussion, Conclusions, and\n# Future Perspectives\") plus a following blank paragraph, right before the\n# existing first paragraph (\"Note: some generic words ...\").\n#\n# We build the new content with raw OOXML (flat-OPC wrapped) and insert it\n# via Range.InsertXML so the run boundaries match the source document\n# exactly (the title text is split across several runs - most likely left\n# over from Word's auto-capitalize-first-letter-of-sentence typing history:\n# \"General Discussion, \" | \"C\" | \"onclusions, and \" | \"F\" | \"uture \" | \"P\" |\n# \"erspectives\"), and the spacer paragraph is a fully empty <w:p/> (no run).\n\n$d = $word.ActiveDocument\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">General Discussion, </w:t></w:r>' +\n  '<w:r><w:t>C</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">onclusions, and </w:t></w:r>' +\n  '<w:r><w:t>F</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">uture </w:t></w:r>' +\n  '<w:r><w:t>P</w:t></w:r>' +\n  '<w:r><w:t>erspectives</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n$r = $d.Range(0, 0)\n$null = $r.InsertXML($flatOpcXml)\n"}
